$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# --- Update the "Status" text everywhere it appears (Ready for handoff -> Handed back) ---
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
}

# --- Update the placeholder handback datetime ---
# zh-cn sheet: datetime becomes 2016-03-23 04:52:14 (re-uses same cells H2/H3)
$wsZh.Cells.Replace("0001-01-01 00:00:00", "2016-03-23 04:52:14")
# de-de sheet: datetime becomes 2016-03-23 04:52:23 (re-uses same cells H2/H3)
$wsDe.Cells.Replace("0001-01-01 00:00:00", "2016-03-23 04:52:23")

# --- Helper to add a "Latest Target File" (F) / "Latest Handback File" (G) pair of
#     hyperlinked cells that mirror the look of the existing hyperlink cells (A / D) ---
function Add-HandbackLinks($ws, $row, $mdUrl, $mdDisplay, $xlfUrl, $xlfDisplay) {
    $fCell = "F" + $row
    $gCell = "G" + $row

    $ws.Hyperlinks.Add($ws.Range($fCell), $mdUrl, "", "", $mdDisplay) | Out-Null
    $ws.Hyperlinks.Add($ws.Range($gCell), $xlfUrl, "", "", $xlfDisplay) | Out-Null

    $ws.Range($fCell).Font.Underline = 2
    $ws.Range($fCell).Font.Color = 15570276
    $ws.Range($gCell).Font.Underline = 2
    $ws.Range($gCell).Font.Color = 15570276
}

# zh-cn sheet (table1) rows 2 & 3 - reuse the same target addresses as columns A / D
Add-HandbackLinks $wsZh 2 `
    "https://github.com/OpenLocalizationTest/oltest/blob/8f87bd8d9a17882f400e62ad533ed0aa84613965/e2e/64799798-104e-43c3-9e33-0204fbd1d975.md" `
    "64799798-104e-43c3-9e33-0204fbd1d975.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/37cc0bc85c84a4580629ee2e66d223354e976209/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/64799798-104e-43c3-9e33-0204fbd1d975.ff01ee6f8dc539f2d8dac0c024eb74c9c7387d87.zh-cn.xlf" `
    "64799798-104e-43c3-9e33-0204fbd1d975.ff01ee6f8dc539f2d8dac0c024eb74c9c7387d87.zh-cn.xlf"

Add-HandbackLinks $wsZh 3 `
    "https://github.com/OpenLocalizationTest/oltest/blob/8f87bd8d9a17882f400e62ad533ed0aa84613965/e2e/6f99df9d-9892-489a-8628-d4551d15ec26.md" `
    "6f99df9d-9892-489a-8628-d4551d15ec26.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/37cc0bc85c84a4580629ee2e66d223354e976209/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6f99df9d-9892-489a-8628-d4551d15ec26.643a91f9656d9b46599fd83514e8cbf501731551.zh-cn.xlf" `
    "6f99df9d-9892-489a-8628-d4551d15ec26.643a91f9656d9b46599fd83514e8cbf501731551.zh-cn.xlf"

# de-de sheet (table2) rows 2 & 3 - reuse the same target addresses as columns A / D
Add-HandbackLinks $wsDe 2 `
    "https://github.com/OpenLocalizationTest/oltest/blob/8f87bd8d9a17882f400e62ad533ed0aa84613965/e2e/64799798-104e-43c3-9e33-0204fbd1d975.md" `
    "64799798-104e-43c3-9e33-0204fbd1d975.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a87943f731586929d2c0d3b2b049d30b4da028d0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/64799798-104e-43c3-9e33-0204fbd1d975.ff01ee6f8dc539f2d8dac0c024eb74c9c7387d87.de-de.xlf" `
    "64799798-104e-43c3-9e33-0204fbd1d975.ff01ee6f8dc539f2d8dac0c024eb74c9c7387d87.de-de.xlf"

Add-HandbackLinks $wsDe 3 `
    "https://github.com/OpenLocalizationTest/oltest/blob/8f87bd8d9a17882f400e62ad533ed0aa84613965/e2e/6f99df9d-9892-489a-8628-d4551d15ec26.md" `
    "6f99df9d-9892-489a-8628-d4551d15ec26.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a87943f731586929d2c0d3b2b049d30b4da028d0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6f99df9d-9892-489a-8628-d4551d15ec26.643a91f9656d9b46599fd83514e8cbf501731551.de-de.xlf" `
    "6f99df9d-9892-489a-8628-d4551d15ec26.643a91f9656d9b46599fd83514e8cbf501731551.de-de.xlf"
